# Updates cryptos list values per the Mar 18 2023 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''27.389.34'
$ws.Range("E2").Value = '  +3.24%  '

# Row 3
$ws.Range("D3").Value = '''1.794.79'
$ws.Range("E3").Value = '  +4.27%  '

# Row 4
$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  +0.23%  '

# Row 5
$ws.Range("D5").Value = '''336.30'
$ws.Range("E5").Value = '  +0.95%  '

# Row 6
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  +0.23%  '

# Row 7
$ws.Range("D7").Value = '''0.3783'
$ws.Range("E7").Value = '  +1.93%  '

# Row 8
$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").Value = '''48.86'
$ws.Range("E8").Value = '  +0.88%  '

# Row 9
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '''0.3432'
$ws.Range("E9").Value = '  +2.28%  '

# Row 10
$ws.Range("D10").Value = '''1.203'
$ws.Range("E10").Value = '  +1.57%  '

# Row 11
$ws.Range("D11").Value = '''0.07497'
$ws.Range("E11").Value = '  +1.36%  '

# Row 12
$ws.Range("D12").Value = '''1.002'
$ws.Range("E12").Value = '  +0.14%  '

# Row 13
$ws.Range("D13").Value = '''22.03'
$ws.Range("E13").Value = '  +9.78%  '

# Row 14
$ws.Range("D14").Value = '''6.493'
$ws.Range("E14").Value = '  +1.84%  '

# Row 15
$ws.Range("D15").Value = '''1.792.93'
$ws.Range("E15").Value = '  +4.20%  '

# Row 16
$ws.Range("D16").Value = '''7.031'
$ws.Range("E16").Value = '  -0.19%  '

# Row 17
$ws.Range("D17").Value = '''0.00001094'
$ws.Range("E17").Value = '  +2.31%  '

# Row 18
$ws.Range("D18").Value = '''0.06627'
$ws.Range("E18").Value = '  +0.05%  '

# Row 19
$ws.Range("D19").Value = '''84.69'
$ws.Range("E19").Value = '  +3.37%  '

# Row 20
$ws.Range("D20").Value = '''1.000'
$ws.Range("E20").Value = '  +0.13%  '

# Row 21
$ws.Range("D21").Value = '''17.36'
$ws.Range("E21").Value = '  +5.00%  '

# Row 22
$ws.Range("D22").Value = '''6.460'
$ws.Range("E22").Value = '  +5.27%  '

# Row 23
$ws.Range("D23").Value = '''27.379.37'
$ws.Range("E23").Value = '  +3.15%  '

# Row 24
$ws.Range("D24").Value = '''12.51'
$ws.Range("E24").Value = '  -1.98%  '

# Row 25
$ws.Range("D25").Value = '''2.464'
$ws.Range("E25").Value = '  +1.16%  '

# Row 26
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = '''2.577'
$ws.Range("E26").Value = '  +7.59%  '

# Row 27
$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D27").Value = '''1.503'
$ws.Range("E27").Value = '  +6.26%  '

# Row 28
$ws.Range("D28").Value = '''21.45'
$ws.Range("E28").Value = '  +10.50%  '

# Row 29
$ws.Range("D29").Value = '''150.76'
$ws.Range("E29").Value = '  -0.55%  '

# Row 30
$ws.Range("D30").Value = '''1.996.48'
$ws.Range("E30").Value = '  +4.50%  '

# Row 31
$ws.Range("D31").Value = '''133.52'
$ws.Range("E31").Value = '  +1.54%  '

# Row 32
$ws.Range("D32").Value = '''4.061'
$ws.Range("E32").Value = '  -1.33%  '

# Row 33
$ws.Range("D33").Value = '''6.137'
$ws.Range("E33").Value = '  +3.15%  '

# Row 34
$ws.Range("D34").Value = '''0.08694'
$ws.Range("E34").Value = '  +1.00%  '

# Row 35
$ws.Range("D35").Value = '''13.30'
$ws.Range("E35").Value = '  +4.62%  '

# Row 36
$ws.Range("D36").Value = '''1.673'
$ws.Range("E36").Value = '  -1.23%  '

# Row 37
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '''5.447'
$ws.Range("E37").Value = '  +1.79%  '

# Row 38
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").Value = '''0.6895'
$ws.Range("E38").Value = '  +11.48%  '

# Row 39
$ws.Range("D39").Value = '''0.06358'
$ws.Range("E39").Value = '  +2.37%  '

# Row 40
$ws.Range("D40").Value = '''0.2201'
$ws.Range("E40").Value = '  +2.36%  '

# Row 41
$ws.Range("D41").Value = '''8.819'
$ws.Range("E41").Value = '  +4.76%  '

# Row 42
$ws.Range("D42").Value = '''0.02344'
$ws.Range("E42").Value = '  +0.89%  '

# Row 43
$ws.Range("D43").Value = '''1.273'
$ws.Range("E43").Value = '  +4.42%  '

# Row 44
$ws.Range("D44").Value = '''14.43'
$ws.Range("E44").Value = '  +1.36%  '

# Row 45
$ws.Range("D45").Value = '''0.6472'
$ws.Range("E45").Value = '  +7.91%  '

# Row 46
$ws.Range("D46").Value = '''1.001'
$ws.Range("E46").Value = '  +0.20%  '

# Row 47
$ws.Range("D47").Value = '''3.852'
$ws.Range("E47").Value = '  -1.45%  '

# Row 48
$ws.Range("D48").Value = '''2.123'
$ws.Range("E48").Value = '  +4.09%  '

# Row 49
$ws.Range("D49").Value = '''129.88'
$ws.Range("E49").Value = '  +0.55%  '

# Row 50
$ws.Range("E50").Value = '  +0.43%  '

# Row 51
$ws.Range("D51").Value = '''79.49'
$ws.Range("E51").Value = '  +3.41%  '
